$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the second data table (CARBON and SILICON) values
$ws.Range("B20").Value = 88.3
$ws.Range("B21").Value = 95.1
$ws.Range("B22").Value = 94.9
$ws.Range("B23").Value = 96.8

# Update the view: scroll so row 2 is the top visible row, and select A13
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A13").Select()
